# Fix the mislabeled 2050 header (which was erroneously left as the
# leftover numeric value from a formula) and remove the "Total" row
# that was appended to each table.

$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3, 5 (index 1,2,3,5): header label in E1 should read "2050"
foreach ($idx in 1, 2, 3, 5) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("E1").Value = "2050"
}

# Sheet 4: header label in E1 should read "2041-2050" (it uses year ranges)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E1").Value = "2041-2050"

# Sheets 1-4: delete row 13, which held the "Total" row
foreach ($idx in 1, 2, 3, 4) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Rows.Item(13).Delete()
}

# Sheet 6: delete row 4, which held the "Total" row
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
